# Applies the two substantive changes from the target commit:
#
# 1. Slide 16's table (shape 3, graphicFrame "Google Shape;213;p29") switches
#    its table style from {58775FA0-062F-48D4-AB20-6FB47848424E} to
#    {F461EA05-D88A-49E8-BC5A-0B53E1CAA387}. PowerPoint requires a method
#    call (ApplyStyle) rather than a plain property assignment for this.
#
# 2. The presentation's theme (currently "Integral") reverts to the
#    built-in default "Office Theme" colour palette. The theme's 12
#    ThemeColorScheme slots are repointed at the stock Office Theme RGB
#    values (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{F461EA05-D88A-49E8-BC5A-0B53E1CAA387}")

# --- 2. Theme colours back to the stock "Office Theme" palette ------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Item(1).RGB  = 0         # dk1      000000
$themeColors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388   # dk2      44546A
$themeColors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407     # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Item(10).RGB = 4697456   # accent6  70AD47
$themeColors.Item(11).RGB = 12673797  # hlink    0563C1
$themeColors.Item(12).RGB = 7491477   # folHlink 954F72
